# Updated description of state machine and counter exercises
#
# 1) Remove the last slide (slide 10 - the FSM / Baud rate counter /
#    Bit counter / Shift register exercise diagram).
# 2) Refresh the cached "today" date field shown on the slide master,
#    every slide layout, and the notes master (05/09/2022 -> 16/09/2022,
#    en-US short form 9/5/22 -> 9/16/22).

$p = $ppt.ActivePresentation

# --- 1. Delete the trailing slide -----------------------------------
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()

# --- 2. Refresh the date placeholders --------------------------------
function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide master (en-US short date form)
Update-DatePlaceholder $p.SlideMaster.Shapes "9/16/22"

# Every slide layout off the master (en-US short date form)
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j).Shapes "9/16/22"
}

# Notes master (en-NO long date form)
Update-DatePlaceholder $p.NotesMaster.Shapes "16/09/2022"
